$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 9457
$ws.Cells.Item(2, 5).Value = 8383
$ws.Cells.Item(2, 6).Value = 0.8864333298086073
$ws.Cells.Item(2, 7).Value = 0.8844692973201097
$ws.Cells.Item(2, 8).Value = 0.09581466306536589
$ws.Cells.Item(2, 9).Value = 0.08474512771438723
$ws.Cells.Item(2, 10).Value = 40726470.77878331
$ws.Cells.Item(2, 11).Value = 14166871.48193765
$ws.Cells.Item(2, 13).Value = 14166871.48193765
$ws.Cells.Item(2, 14).Value = 54893342.26072096
$ws.Cells.Item(2, 15).Value = 800181092.3172001
$ws.Cells.Item(2, 16).Value = 782481285.3132
$ws.Cells.Item(2, 17).Value = 0.0177045816477775
$ws.Cells.Item(2, 18).Value = 0.01810506110221811

# Row 3
$ws.Cells.Item(3, 4).Value = 9643
$ws.Cells.Item(3, 5).Value = 8567
$ws.Cells.Item(3, 6).Value = 0.8884164679041792
$ws.Cells.Item(3, 7).Value = 0.8863025036209394
$ws.Cells.Item(3, 8).Value = 0.09424911249622747
$ws.Cells.Item(3, 9).Value = 0.08353322436945795
$ws.Cells.Item(3, 10).Value = 42515722.80644882
$ws.Cells.Item(3, 11).Value = 14791534.53323031
$ws.Cells.Item(3, 13).Value = 14791534.53323031
$ws.Cells.Item(3, 14).Value = 57307257.33967912
$ws.Cells.Item(3, 15).Value = 837382254.356528
$ws.Cells.Item(3, 16).Value = 819902078.332458
$ws.Cells.Item(3, 17).Value = 0.01766401718722426
$ws.Cells.Item(3, 18).Value = 0.01804061109750299

# Row 4
$ws.Cells.Item(4, 4).Value = 9835
$ws.Cells.Item(4, 5).Value = 8715
$ws.Cells.Item(4, 6).Value = 0.8861209964412812
$ws.Cells.Item(4, 7).Value = 0.8840535605599513
$ws.Cells.Item(4, 8).Value = 0.09311208234732114
$ws.Cells.Item(4, 9).Value = 0.08231606793030065
$ws.Cells.Item(4, 10).Value = 44356356.26072727
$ws.Cells.Item(4, 11).Value = 15391881.79541372
$ws.Cells.Item(4, 13).Value = 15391881.79541372
$ws.Cells.Item(4, 14).Value = 59748238.05614099
$ws.Cells.Item(4, 15).Value = 875322254.7530119
$ws.Cells.Item(4, 16).Value = 857873306.7470582
$ws.Cells.Item(4, 17).Value = 0.01758424592981109
$ws.Cells.Item(4, 18).Value = 0.01794190549392158

# Row 5
$ws.Cells.Item(5, 4).Value = 10024
$ws.Cells.Item(5, 5).Value = 8891
$ws.Cells.Item(5, 6).Value = 0.8869712689545092
$ws.Cells.Item(5, 7).Value = 0.8843246469067038
$ws.Cells.Item(5, 8).Value = 0.09183918725606031
$ws.Cells.Item(5, 9).Value = 0.08121565684241421
$ws.Cells.Item(5, 10).Value = 46281693.98751035
$ws.Cells.Item(5, 11).Value = 16021428.19629553
$ws.Cells.Item(5, 13).Value = 16021428.19629553
$ws.Cells.Item(5, 14).Value = 62303122.18380587
$ws.Cells.Item(5, 15).Value = 913377479.2056578
$ws.Cells.Item(5, 16).Value = 895891373.7426846
$ws.Cells.Item(5, 17).Value = 0.01754086186822668
$ws.Cells.Item(5, 18).Value = 0.01788322632169596

# Row 6
$ws.Cells.Item(6, 4).Value = 10228
$ws.Cells.Item(6, 5).Value = 9091
$ws.Cells.Item(6, 6).Value = 0.8888345717637857
$ws.Cells.Item(6, 7).Value = 0.8865808464989272
$ws.Cells.Item(6, 8).Value = 0.0905776771966571
$ws.Cells.Item(6, 9).Value = 0.08030443372291884
$ws.Cells.Item(6, 10).Value = 48451138.13060883
$ws.Cells.Item(6, 11).Value = 16735055.39106738
$ws.Cells.Item(6, 13).Value = 16735055.39106738
$ws.Cells.Item(6, 14).Value = 65186193.5216762
$ws.Cells.Item(6, 15).Value = 954116855.7025089
$ws.Cells.Item(6, 16).Value = 936525029.8293273
$ws.Cells.Item(6, 17).Value = 0.01753983832383455
$ws.Cells.Item(6, 18).Value = 0.01786930926354118

# Row 7
$ws.Cells.Item(7, 4).Value = 9456
$ws.Cells.Item(7, 5).Value = 8390
$ws.Cells.Item(7, 6).Value = 0.8872673434856176
$ws.Cells.Item(7, 7).Value = 0.8852078497573328
$ws.Cells.Item(7, 8).Value = 0.09666309586593051
$ws.Cells.Item(7, 9).Value = 0.08556693124236725
$ws.Cells.Item(7, 10).Value = 41213626.84427914
$ws.Cells.Item(7, 11).Value = 14411981.36879557
$ws.Cells.Item(7, 13).Value = 14411981.36879557
$ws.Cells.Item(7, 14).Value = 55625608.21307472
$ws.Cells.Item(7, 15).Value = 800758175.3072001
$ws.Cells.Item(7, 16).Value = 783058368.3032
$ws.Cells.Item(7, 17).Value = 0.01799791973808648
$ws.Cells.Item(7, 18).Value = 0.01840473450277369

# Row 8
$ws.Cells.Item(8, 4).Value = 9643
$ws.Cells.Item(8, 5).Value = 8553
$ws.Cells.Item(8, 6).Value = 0.886964637560925
$ws.Cells.Item(8, 7).Value = 0.8848541278708877
$ws.Cells.Item(8, 8).Value = 0.09541183041090721
$ws.Cells.Item(8, 9).Value = 0.08442555198680833
$ws.Cells.Item(8, 10).Value = 43058552.09083918
$ws.Cells.Item(8, 11).Value = 15064526.98515879
$ws.Cells.Item(8, 13).Value = 15064526.98515879
$ws.Cells.Item(8, 14).Value = 58123079.07599795
$ws.Cells.Item(8, 15).Value = 837166427.141628
$ws.Cells.Item(8, 16).Value = 819686251.117558
$ws.Cells.Item(8, 17).Value = 0.0179946621086971
$ws.Cells.Item(8, 18).Value = 0.01837840632878721

# Row 9
$ws.Cells.Item(9, 4).Value = 9833
$ws.Cells.Item(9, 5).Value = 8760
$ws.Cells.Item(9, 6).Value = 0.8908776568697244
$ws.Cells.Item(9, 7).Value = 0.888618381010347
$ws.Cells.Item(9, 8).Value = 0.09408327412763762
$ws.Cells.Item(9, 9).Value = 0.08360412673545399
$ws.Cells.Item(9, 10).Value = 45171746.23605794
$ws.Cells.Item(9, 11).Value = 15801201.92710435
$ws.Cells.Item(9, 13).Value = 15801201.92710435
$ws.Cells.Item(9, 14).Value = 60972948.16316229
$ws.Cells.Item(9, 15).Value = 875679138.500765
$ws.Cells.Item(9, 16).Value = 858230190.4948111
$ws.Cells.Item(9, 17).Value = 0.0180445110913083
$ws.Cells.Item(9, 18).Value = 0.01841137972318848

# Row 10
$ws.Cells.Item(10, 4).Value = 10029
$ws.Cells.Item(10, 5).Value = 8920
$ws.Cells.Item(10, 6).Value = 0.8894206800279191
$ws.Cells.Item(10, 7).Value = 0.8872090710165108
$ws.Cells.Item(10, 8).Value = 0.09311928308151099
$ws.Cells.Item(10, 9).Value = 0.08261627263647085
$ws.Cells.Item(10, 10).Value = 47228158.87806591
$ws.Cells.Item(10, 11).Value = 16496334.53991937
$ws.Cells.Item(10, 13).Value = 16496334.53991937
$ws.Cells.Item(10, 14).Value = 63724493.41798528
$ws.Cells.Item(10, 15).Value = 914254174.7405434
$ws.Cells.Item(10, 16).Value = 896768069.2775702
$ws.Cells.Item(10, 17).Value = 0.01804348833802248
$ws.Cells.Item(10, 18).Value = 0.01839531881772808

# Row 11
$ws.Cells.Item(11, 4).Value = 10232
$ws.Cells.Item(11, 5).Value = 9105
$ws.Cells.Item(11, 6).Value = 0.8898553557466771
$ws.Cells.Item(11, 7).Value = 0.8879461673493271
$ws.Cells.Item(11, 8).Value = 0.09199977791416801
$ws.Cells.Item(11, 9).Value = 0.08169085019587474
$ws.Cells.Item(11, 10).Value = 49468210.12644157
$ws.Cells.Item(11, 11).Value = 17245315.50428019
$ws.Cells.Item(11, 13).Value = 17245315.50428019
$ws.Cells.Item(11, 14).Value = 66713525.63072176
$ws.Cells.Item(11, 15).Value = 955710133.1542411
$ws.Cells.Item(11, 16).Value = 938118307.2810595
$ws.Cells.Item(11, 17).Value = 0.01804450419225281
$ws.Cells.Item(11, 18).Value = 0.01838287918531528

# Row 12
$ws.Cells.Item(12, 4).Value = 9463
$ws.Cells.Item(12, 5).Value = 8393
$ws.Cells.Item(12, 6).Value = 0.8869280355067104
$ws.Cells.Item(12, 7).Value = 0.8855243722304283
$ws.Cells.Item(12, 8).Value = 0.09671184592974379
$ws.Cells.Item(12, 9).Value = 0.08564069665418225
$ws.Cells.Item(12, 10).Value = 41239758.28268903
$ws.Cells.Item(12, 11).Value = 14425047.08800052
$ws.Cells.Item(12, 13).Value = 14425047.08800052
$ws.Cells.Item(12, 14).Value = 55664805.37068957
$ws.Cells.Item(12, 15).Value = 800750889.2872001
$ws.Cells.Item(12, 16).Value = 783051082.2832
$ws.Cells.Item(12, 17).Value = 0.01801440033471731
$ws.Cells.Item(12, 18).Value = 0.0184215914061958

# Row 13
$ws.Cells.Item(13, 4).Value = 9652
$ws.Cells.Item(13, 5).Value = 8568
$ws.Cells.Item(13, 6).Value = 0.8876916701201824
$ws.Cells.Item(13, 7).Value = 0.8864059590316573
$ws.Cells.Item(13, 8).Value = 0.1019541322918984
$ws.Cells.Item(13, 9).Value = 0.09037275041144063
$ws.Cells.Item(13, 10).Value = 48077747.45953142
$ws.Cells.Item(13, 11).Value = 17574124.66950491
$ws.Cells.Item(13, 13).Value = 17574124.66950491
$ws.Cells.Item(13, 14).Value = 65651872.12903633
$ws.Cells.Item(13, 15).Value = 836946297.5090281
$ws.Cells.Item(13, 16).Value = 819466121.4849579
$ws.Cells.Item(13, 17).Value = 0.02099791195899919
$ws.Cells.Item(13, 18).Value = 0.02144582211362047

# Row 14
$ws.Cells.Item(14, 4).Value = 9836
$ws.Cells.Item(14, 5).Value = 8719
$ws.Cells.Item(14, 6).Value = 0.8864375762505083
$ws.Cells.Item(14, 7).Value = 0.8844593223777643
$ws.Cells.Item(14, 8).Value = 0.1063217547815781
$ws.Cells.Item(14, 9).Value = 0.09403726718812941
$ws.Cells.Item(14, 10).Value = 54561128.47940587
$ws.Cells.Item(14, 11).Value = 20495893.04877832
$ws.Cells.Item(14, 13).Value = 20495893.04877832
$ws.Cells.Item(14, 14).Value = 75057021.52818419
$ws.Cells.Item(14, 15).Value = 875020684.362587
$ws.Cells.Item(14, 16).Value = 857571736.3566331
$ws.Cells.Item(14, 17).Value = 0.02342332405971483
$ws.Cells.Item(14, 18).Value = 0.02389991668318558

# Row 15
$ws.Cells.Item(15, 4).Value = 10032
$ws.Cells.Item(15, 5).Value = 8902
$ws.Cells.Item(15, 6).Value = 0.8873604465709729
$ws.Cells.Item(15, 7).Value = 0.8854187388104238
$ws.Cells.Item(15, 8).Value = 0.1094883113450708
$ws.Cells.Item(15, 9).Value = 0.09694300254563555
$ws.Cells.Item(15, 10).Value = 60667658.19851614
$ws.Cells.Item(15, 11).Value = 23216084.20014448
$ws.Cells.Item(15, 13).Value = 23216084.20014448
$ws.Cells.Item(15, 14).Value = 83883742.3986606
$ws.Cells.Item(15, 15).Value = 914085108.69052
$ws.Cells.Item(15, 16).Value = 896599003.2275469
$ws.Cells.Item(15, 17).Value = 0.02539816476542634
$ws.Cells.Item(15, 18).Value = 0.02589349766904937

# Row 16
$ws.Cells.Item(16, 4).Value = 10236
$ws.Cells.Item(16, 5).Value = 9104
$ws.Cells.Item(16, 6).Value = 0.889409925752247
$ws.Cells.Item(16, 7).Value = 0.8878486444314414
$ws.Cells.Item(16, 8).Value = 0.1086977081952903
$ws.Cells.Item(16, 9).Value = 0.09650711287399286
$ws.Cells.Item(16, 10).Value = 64007671.86784674
$ws.Cells.Item(16, 11).Value = 24515046.37498279
$ws.Cells.Item(16, 13).Value = 24515046.37498279
$ws.Cells.Item(16, 14).Value = 88522718.24282953
$ws.Cells.Item(16, 15).Value = 955327879.892617
$ws.Cells.Item(16, 16).Value = 937736054.0194355
$ws.Cells.Item(16, 17).Value = 0.02566139530831905
$ws.Cells.Item(16, 18).Value = 0.02614280027935738
